$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

# Text columns (A-D): briefly force text formatting so date/time-looking
# strings and the zero-padded week number aren't reinterpreted by Excel's
# automatic type inference (which would turn "2025-01-22" into a date
# serial and "03" into the number 3), then drop the style override so the
# new row renders unstyled, like the rest of the data rows.
$c = $ws.Cells.Item($row, 1)
$c.NumberFormat = "@"
$c.Value = "2025-01-22"
$c.Style = "Normal"

$c = $ws.Cells.Item($row, 2)
$c.NumberFormat = "@"
$c.Value = "21:59:02"
$c.Style = "Normal"

$c = $ws.Cells.Item($row, 3)
$c.Value = "Wednesday"

$c = $ws.Cells.Item($row, 4)
$c.NumberFormat = "@"
$c.Value = "03"
$c.Style = "Normal"

# Numeric columns (E-T)
$ws.Cells.Item($row, 5).Value = 126306
$ws.Cells.Item($row, 6).Value = 142125
$ws.Cells.Item($row, 7).Value = 168628
$ws.Cells.Item($row, 8).Value = 158622
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142972
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192294
$ws.Cells.Item($row, 14).Value = 115683
$ws.Cells.Item($row, 15).Value = 45627
$ws.Cells.Item($row, 16).Value = 28447
$ws.Cells.Item($row, 17).Value = 65696
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48227
$ws.Cells.Item($row, 20).Value = -1
